{"js": "// Apply the three textual changes described by the diff (a new \"Quantity\"/\n// \"TotalPrice\" column was wired into Products_ordered, \"Order\" was renamed\n// to \"Orders\", and \"EmailAddress\" was dropped from \"User\"):\n//\n// 1. User(..., Address, EmailAddress)  -> User(..., Address)\n// 2. Order(ID, Total, Username)        -> Orders(ID, Total, Username)\n// 3. Products_ordered(ProductName, OrderID)\n//      -> Products_ordered(ProductName, OrderID, Quantity, TotalPrice)\n//\n// Each edit is scoped to the smallest possible search hit so the\n// surrounding runs (and their underline formatting on the primary-key\n// fields) are left untouched.\n\nconst body = context.document.body;\n\n// --- 1. Drop \", EmailAddress\" from the User(...) relation -------------\nconst emailResults = body.search(\", EmailAddress)\", { matchCase: true });\nemailResults.load(\"items\");\nawait context.sync();\nfor (const r of emailResults.items) {\n  r.insertText(\")\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2. Rename Order(...) to Orders(...) -------------------------------\nconst orderResults = body.search(\"Order(\", { matchCase: true });\norderResults.load(\"items\");\nawait context.sync();\nfor (const r of orderResults.items) {\n  r.insertText(\"Orders(\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 3. Add \", Quantity, TotalPrice\" to Products_ordered(...) ----------\n// Locate the \"OrderID\" occurrence that lives inside the Products_ordered(...)\n// paragraph (the string \"OrderID\" is not unique document-wide, so\n// disambiguate using the owning paragraph's text).\nconst orderIdResults = body.search(\"OrderID\", { matchCase: true });\norderIdResults.load(\"items\");\nawait context.sync();\n\nlet orderIdInProducts = null;\nfor (const r of orderIdResults.items) {\n  const paras = r.paragraphs;\n  paras.load(\"items/text\");\n  await context.sync();\n  if (paras.items.length > 0 && paras.items[0].text.indexOf(\"Products_ordered(\") === 0) {\n    orderIdInProducts = r;\n    break;\n  }\n}\n\nif (orderIdInProducts) {\n  // First extend the underlined key list with \", \" (typing right after\n  // \"OrderID\" naturally continues the underlined run, matching the key\n  // list's own internal separator formatting).\n  const afterComma = orderIdInProducts.insertText(\", \", Word.InsertLocation.after);\n  await context.sync();\n\n  // Then append the two new (non-key / non-underlined) attributes.\n  const newAttrs = afterComma.insertText(\"Quantity, TotalPrice\", Word.InsertLocation.after);\n  await context.sync();\n  newAttrs.font.underline = Word.UnderlineType.none;\n  await context.sync();\n}\n", "ps1": "# Apply the three textual changes described by the diff (a new\n# \"Quantity\"/\"TotalPrice\" column was wired into Products_ordered,\n# \"Order\" was renamed to \"Orders\", and \"EmailAddress\" was dropped\n# from \"User\"):\n#\n# 1. User(..., Address, EmailAddress)  -> User(..., Address)\n# 2. Order(ID, Total, Username)        -> Orders(ID, Total, Username)\n# 3. Products_ordered(ProductName, OrderID)\n#      -> Products_ordered(ProductName, OrderID, Quantity, TotalPrice)\n#\n# Each edit is scoped to the owning paragraph (found by its distinctive\n# leading text) so the surrounding runs -- and the underline formatting\n# on the primary-key fields -- are left untouched.\n\n$d = $word.ActiveDocument\n\n# --- 1. Drop \", EmailAddress\" from the User(...) relation --------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"User(\")) {\n        $r = $p.Range\n        # wdFindContinue = 1, wdReplaceOne = 2\n        $null = $r.Find.Execute(\", EmailAddress)\", $false, $false, $false, $false, $false, $true, 1, $false, \")\", 2)\n        break\n    }\n}\n\n# --- 2. Rename Order(...) to Orders(...) --------------------------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Order(\")) {\n        $r = $p.Range\n        $null = $r.Find.Execute(\"Order(\", $false, $false, $false, $false, $false, $true, 1, $false, \"Orders(\", 2)\n        break\n    }\n}\n\n# --- 3. Add \", Quantity, TotalPrice\" to Products_ordered(...) ----------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Products_ordered(\")) {\n        $r = $p.Range\n        $null = $r.Find.Execute(\"OrderID\")\n\n        # Collapse to a point right after \"OrderID\" and type \", \" --\n        # inserting here naturally continues the underlined run, matching\n        # the key list's own separator formatting.\n        $r.Collapse(0)   # wdCollapseEnd\n        $r.InsertAfter(\", \")\n\n        # Collapse again (now past the \", \") and add the two new,\n        # non-key attributes, then explicitly clear underline on them.\n        $r.Collapse(0)\n        $r.InsertAfter(\"Quantity, TotalPrice\")\n        $r.Font.Underline = 0   # wdUnderlineNone\n        break\n    }\n}\n"}
